$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.010.99'
$ws.Range('E2').Value = '  -3.17%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.797.44'
$ws.Range('E3').Value = '  -3.41%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  +0.03%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '307.14'
$ws.Range('E5').Value = '  -3.02%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  +0.01%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4189'
$ws.Range('E7').Value = '  -3.33%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3574'
$ws.Range('E8').Value = '  -3.95%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07083'
$ws.Range('E9').Value = '  -4.09%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8447'
$ws.Range('E10').Value = '  -4.51%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.15'
$ws.Range('E11').Value = '  -5.17%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.795.42'
$ws.Range('E12').Value = '  -5.15%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.291'
$ws.Range('E13').Value = '  -3.80%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.352'
$ws.Range('E14').Value = '  -4.33%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.06762'
$ws.Range('E15').Value = '  -2.90%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.004'
$ws.Range('E16').Value = '  +0.15%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '79.71'
$ws.Range('E17').Value = '  -2.50%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008689'
$ws.Range('E18').Value = '  -4.64%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.99'
$ws.Range('E20').Value = '  -4.11%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '26.905.72'
$ws.Range('E21').Value = '  -5.02%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.047'
$ws.Range('E22').Value = '  -1.52%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.94'
$ws.Range('E23').Value = '  -1.10%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.992.49'
$ws.Range('E24').Value = '  -6.12%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.932'
$ws.Range('E25').Value = '  -2.80%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '152.70'
$ws.Range('E26').Value = '  -1.56%  '

$ws.Range('E27').Value = '  -6.22%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.009'
$ws.Range('E28').Value = '  -6.71%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '112.90'
$ws.Range('E29').Value = '  -3.18%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.642'
$ws.Range('E30').Value = '  -12.35%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08959'
$ws.Range('E31').Value = '  -0.17%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7197'
$ws.Range('E32').Value = '  -9.31%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.857'
$ws.Range('E33').Value = '  -4.55%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.297'
$ws.Range('E34').Value = '  -7.86%  '

$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.077'
$ws.Range('E35').Value = '  -9.03%  '

$ws.Range('B36').Value = 'Frax'
$ws.Range('C36').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.000'
$ws.Range('E36').Value = '  +0.05%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.073'
$ws.Range('E37').Value = '  -3.70%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01900'
$ws.Range('E38').Value = '  -3.70%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05114'
$ws.Range('E39').Value = '  -6.63%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.1625'
$ws.Range('E40').Value = '  -4.29%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4944'
$ws.Range('E41').Value = '  -5.31%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.581'
$ws.Range('E42').Value = '  -9.43%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.962'
$ws.Range('E43').Value = '  -12.47%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.040'
$ws.Range('E44').Value = '  -8.24%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '104.52'
$ws.Range('E45').Value = '  -2.66%  '

$ws.Range('E46').Value = '  -4.25%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.9999'
$ws.Range('E47').Value = '  +0.03%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.06298'
$ws.Range('E48').Value = '  -4.50%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.4526'
$ws.Range('E49').Value = '  -6.18%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.599'
$ws.Range('E50').Value = '  -5.02%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.692'
$ws.Range('E51').Value = '  -9.61%  '
